$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy the formatting (style, number format, alignment, borders) of column F
# into the new column G before writing values, so the new cells inherit the
# same look (currency format for the data rows, bold/border header style for
# row 1 and the right-aligned currency total style for row 52).
$ws.Range("F1:F52").Copy()
$ws.Range("G1:G52").PasteSpecial(-4122)

# Header
$ws.Range("G1").Value = "PRESUPUESTO"

# Data rows (2-51) hold a budget value of 0
$ws.Range("G2:G51").Value = 0

# Totals row
$ws.Range("G52").Value = 0

# Set the new column's width to match the authored width of 17
$ws.Columns.Item(7).ColumnWidth = 16.17
